# refactor: extract checkpoint functionality into separate module
#
# Replace the report rows (rows 2-7) with the new dataset. The shared
# strings table is regenerated by the engine at save time from the order
# in which string values are assigned, so the cells are deliberately set
# column-by-column (all of column A, then all of column B, etc.) to
# reproduce the canonical ordering of the new shared-string entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Item numbers -------------------------------------------------
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 13
$ws.Range("A4").Value = 14
$ws.Range("A5").Value = 16
$ws.Range("A6").Value = 17
$ws.Range("A7").Value = 18

# --- Column B: Descrição Original ------------------------------------------
$ws.Range("B2").Value = "RISPERIDONA 3MG"
$ws.Range("B3").Value = "RISPERIDONA 1MG"
$ws.Range("B4").Value = "LEVOMEPROMAZINA 4% GOTAS"
$ws.Range("B5").Value = "LIDOCAINA 2% C/ VASO CONSTRITO"
$ws.Range("B6").Value = "NITRATO DE CERIO +SULFADIAZINA"
$ws.Range("B7").Value = "COLAGENASE+CLORAFENICOL POMADA 30g"

# --- Column C: Descrição Final ----------------------------------------------
$ws.Range("C2").Value = "RISPERIDONA"
$ws.Range("C3").Value = "RISPERIDONA"
$ws.Range("C4").Value = "levomepromazina"
$ws.Range("C5").Value = "LIDOCAÍNA"
$ws.Range("C6").Value = "NITRATO DE CERIO;SULFADIAZINA"
$ws.Range("C7").Value = "COLAGENASE"

# --- Column D: Concentração ---------------------------------------------------
# "3mg" / "1mg" / "30g" round-trip fine as plain text, but a bare "4%" or
# "2%" gets auto-recognized as a percentage number by the COM value
# setter. Force those two through a text number-format so they stay
# plain strings, then drop the number format again so the cell ends up
# with no explicit style, matching the rest of the column.
$ws.Range("D2").Value = "3mg"
$ws.Range("D3").Value = "1mg"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4%"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2%"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "Concentração não encontrada"
$ws.Range("D7").Value = "30g"

# --- Column E: Laboratório Original -----------------------------------------
$ws.Range("E2").Value = "Vitamedic"
$ws.Range("E3").Value = "Vitamedic"
$ws.Range("E4").Value = "Sanofi"
$ws.Range("E5").Value = "Cristália"
$ws.Range("E6").Value = "Cristália"
$ws.Range("E7").Value = "Cristália"

# --- Column F: Laboratório Final ---------------------------------------------
$ws.Range("F2").Value = "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA"
$ws.Range("F3").Value = "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA"
$ws.Range("F4").Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Range("F5").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F6").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F7").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
